$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2057416267942584
$ws.Range("C2").Value = 0.5789473684210527
$ws.Range("J2").Value = 0.01913875598086124
$ws.Range("P2").Value = 0.1244019138755981
$ws.Range("S2").Value = 0.07177033492822966
$ws.Range("C3").Value = 0.024
$ws.Range("J3").Value = 0.016
$ws.Range("P3").Value = 0.744
$ws.Range("S3").Value = 0.216
$ws.Range("J4").Value = 0.06060606060606061
$ws.Range("P4").Value = 0.6363636363636364
$ws.Range("S4").Value = 0.303030303030303
$ws.Range("B6").Value = 0.07647058823529412
$ws.Range("D6").Value = 0.01176470588235294
$ws.Range("F6").Value = 0.1
$ws.Range("J6").Value = 0.1882352941176471
$ws.Range("O6").Value = 0.02941176470588235
$ws.Range("Q6").Value = 0.1411764705882353
$ws.Range("R6").Value = 0.07058823529411765
$ws.Range("S6").Value = 0.3823529411764706
$ws.Range("B7").Value = 0.1052631578947368
$ws.Range("D7").Value = 0.06015037593984962
$ws.Range("F7").Value = 0.04511278195488722
$ws.Range("J7").Value = 0.1353383458646616
$ws.Range("O7").Value = 0.02255639097744361
$ws.Range("Q7").Value = 0.1428571428571428
$ws.Range("R7").Value = 0.1278195488721804
$ws.Range("S7").Value = 0.3609022556390977
$ws.Range("B8").Value = 0.07624633431085044
$ws.Range("D8").Value = 0.008797653958944282
$ws.Range("F8").Value = 0.07624633431085044
$ws.Range("J8").Value = 0.1026392961876833
$ws.Range("O8").Value = 0.01173020527859238
$ws.Range("Q8").Value = 0.1994134897360704
$ws.Range("R8").Value = 0.126099706744868
$ws.Range("S8").Value = 0.3988269794721407
$ws.Range("B9").Value = 0.1235294117647059
$ws.Range("D9").Value = 0.01764705882352941
$ws.Range("F9").Value = 0.05882352941176471
$ws.Range("J9").Value = 0.05294117647058823
$ws.Range("O9").Value = 0.02352941176470588
$ws.Range("Q9").Value = 0.1529411764705882
$ws.Range("R9").Value = 0.1588235294117647
$ws.Range("S9").Value = 0.4117647058823529
$ws.Range("B10").Value = 0.08708414872798434
$ws.Range("D10").Value = 0.01663405088062622
$ws.Range("E10").Value = 0.0009784735812133072
$ws.Range("F10").Value = 0.06164383561643835
$ws.Range("J10").Value = 0.1213307240704501
$ws.Range("O10").Value = 0.01761252446183953
$ws.Range("Q10").Value = 0.2289628180039139
$ws.Range("R10").Value = 0.1076320939334638
$ws.Range("S10").Value = 0.3581213307240704
$ws.Range("G11").Value = 0.1477272727272727
$ws.Range("J11").Value = 0.1098484848484848
$ws.Range("K11").Value = 0.2424242424242424
$ws.Range("L11").Value = 0.4810606060606061
$ws.Range("S11").Value = 0.01893939393939394
$ws.Range("G12").Value = 0.6201550387596899
$ws.Range("J12").Value = 0.2713178294573643
$ws.Range("K12").Value = 0.0310077519379845
$ws.Range("L12").Value = 0.0310077519379845
$ws.Range("S12").Value = 0.04651162790697674
$ws.Range("G13").Value = 0.6428571428571429
$ws.Range("J13").Value = 0.3571428571428572
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.00510204081632653
$ws.Range("H15").Value = 0.1275510204081633
$ws.Range("I15").Value = 0.1020408163265306
$ws.Range("J15").Value = 0.3928571428571428
$ws.Range("K15").Value = 0.05102040816326531
$ws.Range("O15").Value = 0.05612244897959184
$ws.Range("S15").Value = 0.2653061224489796
$ws.Range("F16").Value = 0.01492537313432836
$ws.Range("H16").Value = 0.1343283582089552
$ws.Range("I16").Value = 0.06716417910447761
$ws.Range("J16").Value = 0.4328358208955224
$ws.Range("K16").Value = 0.1417910447761194
$ws.Range("M16").Value = 0.04477611940298507
$ws.Range("O16").Value = 0.02985074626865672
$ws.Range("S16").Value = 0.1343283582089552
$ws.Range("F17").Value = 0.008
$ws.Range("H17").Value = 0.1653333333333333
$ws.Range("I17").Value = 0.09866666666666667
$ws.Range("J17").Value = 0.432
$ws.Range("K17").Value = 0.09866666666666667
$ws.Range("M17").Value = 0.01333333333333333
$ws.Range("N17").Value = 0.002666666666666667
$ws.Range("O17").Value = 0.06133333333333333
$ws.Range("S17").Value = 0.12
$ws.Range("F18").Value = 0.01435406698564593
$ws.Range("H18").Value = 0.1770334928229665
$ws.Range("I18").Value = 0.1004784688995215
$ws.Range("J18").Value = 0.4210526315789473
$ws.Range("K18").Value = 0.1052631578947368
$ws.Range("M18").Value = 0.01435406698564593
$ws.Range("O18").Value = 0.05741626794258373
$ws.Range("S18").Value = 0.1100478468899522
$ws.Range("F19").Value = 0.019211324570273
$ws.Range("H19").Value = 0.2042467138523761
$ws.Range("I19").Value = 0.08291203235591507
$ws.Range("J19").Value = 0.3508594539939333
$ws.Range("K19").Value = 0.109201213346815
$ws.Range("M19").Value = 0.01516683518705763
$ws.Range("O19").Value = 0.09605662285136501
$ws.Range("S19").Value = 0.1223458038422649
